$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33: H33(old 358.45947),I33(old 220.5),J33(old 613.1539),K33(old 220.5),L33(old 613.1539),M33(old 8.5),N33(old -1071.1539)
$ws.Range("H33").Value = 337.05
$ws.Range("I33").Value = 209.07692
$ws.Range("J33").Value = 574.7143
$ws.Range("K33").Value = 209.07692
$ws.Range("L33").Value = 574.7143
$ws.Range("M33").Value = 19.92308
$ws.Range("N33").Value = -1032.7143
# Row 76: H76(old 7234.875),I76(old 5937.5),K76(old 5937.5),M76(old -5622.5)
$ws.Range("H76").Value = 6780
$ws.Range("I76").Value = 5715.2
$ws.Range("K76").Value = 5715.2
$ws.Range("M76").Value = -5400.2
# Row 79: H79(old 7234.875),I79(old 5937.5),K79(old 5937.5),M79(old -4845.5)
$ws.Range("H79").Value = 6780
$ws.Range("I79").Value = 5715.2
$ws.Range("K79").Value = 5715.2
$ws.Range("M79").Value = -4623.2
# Row 115: H115(old 561.3570999999999),I115(old 561.3570999999999),K115(old 1684.0713),M115(old -117.0712999999998)
$ws.Range("H115").Value = 536.2
$ws.Range("I115").Value = 536.2
$ws.Range("K115").Value = 1608.6
$ws.Range("M115").Value = -41.60000000000014
# Row 118: H118(old 2767.8),I118(old 3280),K118(old 9840),M118(old -8183)
$ws.Range("H118").Value = 1833.375
$ws.Range("I118").Value = 1778
$ws.Range("K118").Value = 5334
$ws.Range("M118").Value = -3677
# Row 132: H132(old 45643.926),I132(old 48036.46),K132(old 144109.38),M132(old -141579.38)
$ws.Range("H132").Value = 46516.04
$ws.Range("I132").Value = 49010.797
$ws.Range("K132").Value = 147032.391
$ws.Range("M132").Value = -144502.391
# Row 135: H135(old 1793.2632),I135(old 1886.6471),K135(old 16979.8239),M135(old -14444.8239)
$ws.Range("H135").Value = 1891.1111
$ws.Range("I135").Value = 2002.5625
$ws.Range("K135").Value = 18023.0625
$ws.Range("M135").Value = -15488.0625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 38: H38(old 1700),I38(old 1700),K38(old 1700),M38(old -1233)
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
# Row 61: H61(old 2579778.8),I61(old 2579778.8),K61(old 2579778.8),M61(old -2579566.8)
$ws.Range("H61").Value = 2236178.2
$ws.Range("I61").Value = 2236178.2
$ws.Range("K61").Value = 2236178.2
$ws.Range("M61").Value = -2235966.2
# Row 74: H74(old 2018668.2),I74(old 3127133.5),K74(old 3127133.5),M74(old -3126259.5)
$ws.Range("H74").Value = 1925537.9
$ws.Range("I74").Value = 2909020.2
$ws.Range("K74").Value = 2909020.2
$ws.Range("M74").Value = -2908146.2
# Row 77: H77(old 2018668.2),I77(old 3127133.5),K77(old 15635667.5),M77(old -15631299.5)
$ws.Range("H77").Value = 1925537.9
$ws.Range("I77").Value = 2909020.2
$ws.Range("K77").Value = 14545101
$ws.Range("M77").Value = -14540733
# Row 136: H136(old 2579778.8),I136(old 2579778.8),K136(old 7739336.399999999),M136(old -7736786.399999999)
$ws.Range("H136").Value = 2236178.2
$ws.Range("I136").Value = 2236178.2
$ws.Range("K136").Value = 6708534.600000001
$ws.Range("M136").Value = -6705984.600000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 64: H64(old 2153.7856),J64(old 2908.5),L64(old 2908.5),N64(old -3358.5)
$ws.Range("H64").Value = 2296.3845
$ws.Range("J64").Value = 3198.3333
$ws.Range("L64").Value = 3198.3333
$ws.Range("N64").Value = -3648.3333
# Row 67: H67(old 2153.7856),J67(old 2908.5),L67(old 2908.5),N67(old -4468.5)
$ws.Range("H67").Value = 2296.3845
$ws.Range("J67").Value = 3198.3333
$ws.Range("L67").Value = 3198.3333
$ws.Range("N67").Value = -4758.3333
# Row 86: H86(old 3235.5),I86(old 3235.5),K86(old 3235.5),M86(old -2112.5)
$ws.Range("H86").Value = 3139.125
$ws.Range("I86").Value = 3139.125
$ws.Range("K86").Value = 3139.125
$ws.Range("M86").Value = -2016.125
# Row 89: H89(old 3235.5),I89(old 3235.5),K89(old 16177.5),M89(old -10561.5)
$ws.Range("H89").Value = 3139.125
$ws.Range("I89").Value = 3139.125
$ws.Range("K89").Value = 15695.625
$ws.Range("M89").Value = -10079.625
# Row 99: H99(old 88829.664),I99(old 129747.25),J99(old 6994.5),K99(old 129747.25),L99(old 6994.5),M99(old -128249.25),N99(old -9990.5)
$ws.Range("H99").Value = 16744.795
$ws.Range("I99").Value = 19253.115
$ws.Range("J99").Value = 6990.222
$ws.Range("K99").Value = 19253.115
$ws.Range("L99").Value = 6990.222
$ws.Range("M99").Value = -17755.115
$ws.Range("N99").Value = -9986.222
# Row 134: H134(old 830264.0600000001),I134(old 993449.0600000001),K134(old 2980347.18),M134(old -2977812.18)
$ws.Range("H134").Value = 1185904.4
$ws.Range("I134").Value = 1552018.5
$ws.Range("K134").Value = 4656055.5
$ws.Range("M134").Value = -4653520.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16: H16(old 709.6923),I16(old 729.9091),J16(old 598.5),K16(old 729.9091),L16(old 598.5),M16(old -442.9091),N16(old -1172.5)
$ws.Range("H16").Value = 702
$ws.Range("I16").Value = 710.6667
$ws.Range("J16").Value = 598
$ws.Range("K16").Value = 710.6667
$ws.Range("L16").Value = 598
$ws.Range("M16").Value = -423.6667
$ws.Range("N16").Value = -1172
# Row 22: H22(old 1240.4166),J22(old 1482.1666),L22(old 1482.1666),N22(old -2182.1666)
$ws.Range("H22").Value = 1262.5454
$ws.Range("J22").Value = 1579.2
$ws.Range("L22").Value = 1579.2
$ws.Range("N22").Value = -2279.2
# Row 31: H31(old 7256.0967),I31(old 1866.7693),K31(old 1866.7693),M31(old -1571.7693)
$ws.Range("H31").Value = 6870.788
$ws.Range("I31").Value = 1737.6666
$ws.Range("K31").Value = 1737.6666
$ws.Range("M31").Value = -1442.6666
# Row 34: H34(old 7256.0967),I34(old 1866.7693),K34(old 1866.7693),M34(old -1664.7693)
$ws.Range("H34").Value = 6870.788
$ws.Range("I34").Value = 1737.6666
$ws.Range("K34").Value = 1737.6666
$ws.Range("M34").Value = -1535.6666
# Row 105: H105(old 49971.855),I105(old 49971.855),K105(old 49971.855),M105(old -48224.855)
$ws.Range("H105").Value = 49921.57
$ws.Range("I105").Value = 49921.57
$ws.Range("K105").Value = 49921.57
$ws.Range("M105").Value = -48174.57
# Row 113: H113(old 709.6923),I113(old 729.9091),J113(old 598.5),K113(old 729.9091),L113(old 598.5),M113(old 1440.0909),N113(old -4938.5)
$ws.Range("H113").Value = 702
$ws.Range("I113").Value = 710.6667
$ws.Range("J113").Value = 598
$ws.Range("K113").Value = 710.6667
$ws.Range("L113").Value = 598
$ws.Range("M113").Value = 1459.3333
$ws.Range("N113").Value = -4938
# Row 134: H134(old 1227.3784),I134(old 1183.4),K134(old 3550.2),M134(old -1015.2)
$ws.Range("H134").Value = 1179.5641
$ws.Range("I134").Value = 1135.3784
$ws.Range("K134").Value = 3406.1352
$ws.Range("M134").Value = -871.1352000000002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 26: H26(old 216.2)
$ws.Range("H26").Value = 129
# Row 44: H44(old 304.33334),I44(old 499),J44(old 207),K44(old 1497),L44(old 621),M44(old -1099),N44(old -1417)
$ws.Range("H44").Value = 418
$ws.Range("I44").Value = 587
$ws.Range("J44").Value = 192.66667
$ws.Range("K44").Value = 1761
$ws.Range("L44").Value = 578.00001
$ws.Range("M44").Value = -1363
$ws.Range("N44").Value = -1374.00001
# Row 69: H69(old 1067.8334),I69(old 1090.909),J69(old 814),K69(old 3272.727),L69(old 2442),M69(old -2461.727),N69(old -4064)
$ws.Range("H69").Value = 1241.3334
$ws.Range("I69").Value = 1047
$ws.Range("J69").Value = 2504.5
$ws.Range("K69").Value = 3141
$ws.Range("L69").Value = 7513.5
$ws.Range("M69").Value = -2330
$ws.Range("N69").Value = -9135.5
# Row 72: H72(old 1067.8334),I72(old 1090.909),J72(old 814),K72(old 9818.181),L72(old 7326),M72(old -5762.181),N72(old -15438)
$ws.Range("H72").Value = 1241.3334
$ws.Range("I72").Value = 1047
$ws.Range("J72").Value = 2504.5
$ws.Range("K72").Value = 9423
$ws.Range("L72").Value = 22540.5
$ws.Range("M72").Value = -5367
$ws.Range("N72").Value = -30652.5
# Row 138: H138(old 3227.158),I138(old 3227.158),K138(old 9681.474),M138(old -4541.474)
$ws.Range("H138").Value = 3145.6316
$ws.Range("I138").Value = 3145.6316
$ws.Range("K138").Value = 9436.8948
$ws.Range("M138").Value = -4296.8948
# Row 139: H139(old 2990.9),I139(old 1999.2),J139(old 3982.6),K139(old 5997.6),L139(old 11947.8),M139(old -857.6000000000004),N139(old -22227.8)
$ws.Range("H139").Value = 2642.923
$ws.Range("I139").Value = 1799.2858
$ws.Range("J139").Value = 3627.1667
$ws.Range("K139").Value = 5397.857400000001
$ws.Range("L139").Value = 10881.5001
$ws.Range("M139").Value = -257.8574000000008
$ws.Range("N139").Value = -21161.5001
# Row 141: H141(old 5205.4),I141(old 5205.4),K141(old 15616.2),M141(old -10436.2)
$ws.Range("H141").Value = 4887.5
$ws.Range("I141").Value = 4887.5
$ws.Range("K141").Value = 14662.5
$ws.Range("M141").Value = -9482.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46: H46(old 21825),I46(old 6526),J46(old 37124),K46(old 6526),L46(old 37124),M46(old -6370),N46(old -37436)
$ws.Range("H46").Value = 25257.666
$ws.Range("I46").Value = 13349.667
$ws.Range("J46").Value = 37165.668
$ws.Range("K46").Value = 13349.667
$ws.Range("L46").Value = 37165.668
$ws.Range("M46").Value = -13193.667
$ws.Range("N46").Value = -37477.668
# Row 102: H102(old 2128.2068),I102(old 1475.3),K102(old 1475.3),M102(old 146.7)
$ws.Range("H102").Value = 2107.862
$ws.Range("I102").Value = 1445.8
$ws.Range("K102").Value = 1445.8
$ws.Range("M102").Value = 176.2
# Row 126: H126(old 1114977.6),J126(old 4727.273),L126(old 14181.819),N126(old -19121.819)
$ws.Range("H126").Value = 1194261.8
$ws.Range("J126").Value = 4700
$ws.Range("L126").Value = 14100
$ws.Range("N126").Value = -19040
# Row 132: H132(old 525444.0600000001),I132(old 603611.1),K132(old 1810833.3),M132(old -1808303.3)
$ws.Range("H132").Value = 503591.75
$ws.Range("I132").Value = 574914.75
$ws.Range("K132").Value = 1724744.25
$ws.Range("M132").Value = -1722214.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55: H55(old 287.5),J55(old 275.7143),L55(old 275.7143),N55(old -621.7143)
$ws.Range("H55").Value = 297.33334
$ws.Range("J55").Value = 288.5
$ws.Range("L55").Value = 288.5
$ws.Range("N55").Value = -634.5
# Row 82: H82(old 2420.1177),I82(old 2279.4285),J82(old 2518.6),K82(old 2279.4285),L82(old 2518.6),M82(old -1918.4285),N82(old -3240.6)
$ws.Range("H82").Value = 2357.7896
$ws.Range("I82").Value = 2202
$ws.Range("J82").Value = 2471.0908
$ws.Range("K82").Value = 2202
$ws.Range("L82").Value = 2471.0908
$ws.Range("M82").Value = -1841
$ws.Range("N82").Value = -3193.0908
# Row 85: H85(old 2420.1177),I85(old 2279.4285),J85(old 2518.6),K85(old 2279.4285),L85(old 2518.6),M85(old -1031.4285),N85(old -5014.6)
$ws.Range("H85").Value = 2357.7896
$ws.Range("I85").Value = 2202
$ws.Range("J85").Value = 2471.0908
$ws.Range("K85").Value = 2202
$ws.Range("L85").Value = 2471.0908
$ws.Range("M85").Value = -954
$ws.Range("N85").Value = -4967.0908
# Row 136: H136(old 5116.913),I136(old 4804.0454),K136(old 14412.1362),M136(old -11862.1362)
$ws.Range("H136").Value = 4686.6294
$ws.Range("I136").Value = 4405.346
$ws.Range("K136").Value = 13216.038
$ws.Range("M136").Value = -10666.038

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62: H62(old 900000),I62(old 900000),J62(old 0),K62(old 900000),L62(old 0),M62(old -899376),N62(old None)
$ws.Range("H62").Value = 118020
$ws.Range("I62").Value = 450501
$ws.Range("J62").Value = 7193
$ws.Range("K62").Value = 450501
$ws.Range("L62").Value = 7193
$ws.Range("M62").Value = -449877
$ws.Range("N62").Value = -8441
# Row 65: H65(old 900000),I65(old 900000),J65(old 0),K65(old 4500000),L65(old 0),M65(old -4496880),N65(old None)
$ws.Range("H65").Value = 118020
$ws.Range("I65").Value = 450501
$ws.Range("J65").Value = 7193
$ws.Range("K65").Value = 2252505
$ws.Range("L65").Value = 7193
$ws.Range("M65").Value = -2249385
$ws.Range("N65").Value = -42205
# Row 132: H132(old 5034542.5),I132(old 6710760),J132(old 5889.6),K132(old 20132280),L132(old 17668.8),M132(old -20129750),N132(old -22728.8)
$ws.Range("H132").Value = 6496000.5
$ws.Range("I132").Value = 8753005
$ws.Range("J132").Value = 7112.375
$ws.Range("K132").Value = 26259015
$ws.Range("L132").Value = 21337.125
$ws.Range("M132").Value = -26256485
$ws.Range("N132").Value = -26397.125
# Row 136: H136(old 10589856),I136(old 13144753),J136(old 5285),K136(old 39434259),L136(old 15855),M136(old -39431709),N136(old -20955)
$ws.Range("H136").Value = 10032600
$ws.Range("I136").Value = 12706661
$ws.Range("J136").Value = 4872.75
$ws.Range("K136").Value = 38119983
$ws.Range("L136").Value = 14618.25
$ws.Range("M136").Value = -38117433
$ws.Range("N136").Value = -19718.25
# Row 141: H141(old 0),J141(old 0),L141(old 0),N141(old None)
$ws.Range("H141").Value = 95500
$ws.Range("J141").Value = 95500
$ws.Range("L141").Value = 95500
$ws.Range("N141").Value = -105860
